$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (Strikeouts). Update computed values for rows 2-6.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
